# Insert a new weekly price record as row 23 on the "Espinaca" sheet,
# pushing the existing rows 23:84 down to 24:85 (dimension A1:R84 -> A1:R85).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(23).Insert()

$ws.Range("A23").Value = 4
$ws.Range("B23").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C23").Value = "Los Lagos"
$ws.Range("D23").Value = "9/8/2023"
$ws.Range("E23").Value = 10
$ws.Range("F23").Value = 100112012
$ws.Range("G23").Value = "Espinaca"
$ws.Range("H23").Value = "Sin especificar"
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 35
$ws.Range("K23").Value = 13000
$ws.Range("L23").Value = 13000
$ws.Range("M23").Value = 13000
$ws.Range("N23").Value = "$/cuna 10 kilos"
$ws.Range("O23").Value = "Región Metropolitana"
$ws.Range("P23").Value = 1300
$ws.Range("Q23").Value = 10
$ws.Range("R23").Value = "Hortaliza"
